# finish 20 faults in each file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 (fault 16): mutation test column now passes
$ws.Range("F18").Value = 1

# Row 19 (fault 17): coverage test column now fails
$ws.Range("E19").Value = 0

# Row 21 (fault 19): coverage test column now fails
$ws.Range("E21").Value = 0

# Row 22 (fault 20): coverage test column now fails
$ws.Range("E22").Value = 0

# Move selection to reflect where the user ended up after finishing the edits
$ws.Range("F25").Select()
